$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '317.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.91%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.96%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.139'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.51%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08193'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.44%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.028'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.31%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.339'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.37%'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '8.337'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '3.96%'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9415'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.58%'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1349'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-7.85%'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1984'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.07%'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09102'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.16%'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03497'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.19%'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09779'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.03%'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001414'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.35%'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006086'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '4.41%'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.692'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.45%'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.233'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-5.20%'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.38%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1315'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.97%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.962'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.67%'

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.42%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04367'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.22%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001232'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.44%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004799'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '12.33%'

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.04%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004003'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-10.00%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02250'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '10.90%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05186'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.72%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007754'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.07%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009945'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.66%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1399'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.37%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002050'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.57%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009111'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-8.06%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006587'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.09%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.13%'

# Row 48
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.001691'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-6.27%'

# Row 49
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002949'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '2.61%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.13%'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.13%'
